# set list default select to newest
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows of data (distinct event_name / sn_id selector + two urls)
$ws.Cells.Item(18, 1).Value = 43282
$ws.Cells.Item(18, 1).NumberFormat = $ws.Cells.Item(10, 1).NumberFormat
$ws.Cells.Item(18, 2).Value = "預設選單顯示用distinct event_name和sn_id"

$ws.Cells.Item(19, 1).Value = 43283
$ws.Cells.Item(19, 1).NumberFormat = $ws.Cells.Item(10, 1).NumberFormat
$ws.Cells.Item(19, 2).Value = "http://mabitool.tk/SNshare/mabi_get.php"

$ws.Cells.Item(20, 2).Value = "http://mabitool.ddns.net/SNshare/mabi_get.php"

# move the active selection to reflect newest entry being highlighted
$ws.Range("H24").Select()
